# Updated BOM (all components in stock on LCSC on 2/6/2020)
#
# Row 3 of the BOM sheet (component #2, the "309K" 0402 resistor, R2/R4)
# gets re-sourced to a different manufacturer/part:
#   F3 Manufacturer Part : RC-02W3093FT        -> WR04X3093FTL
#   G3 Manufacturer      : Guangdong Fenghua... -> Walsin Tech Corp
#   I3 Supplier Part     : C321438              -> C334683
# (H3 Supplier stays "LCSC"; B3/C3/D3/E3/J3 are unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM_TSDZ2-ESP32-v3")

$ws.Range("F3").Value = "WR04X3093FTL"
$ws.Range("G3").Value = "Walsin Tech Corp"
$ws.Range("I3").Value = "C334683"

# Reflect the author's on-screen state: the whole row 3 was selected
# while making the edit.
$ws.Rows.Item(3).Select()
